$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D values are textual (e.g. "26.515.07" or "1.001") and must
# stay as text, not be auto-converted to numbers by Excel. We force
# the cell to text format before assignment, then restore the default
# "Normal" style so no stray style index is left on the cell.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.515.07"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +6.46%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.721.95"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +3.48%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "332.03"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.18%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.04%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3708"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.84%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "48.15"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.95%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3357"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.65%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.185"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.29%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07408"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.42%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.001"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.09%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.398"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +5.55%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.08"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.69%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.033"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +6.13%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.719.00"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.44%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001069"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.87%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06603"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.37%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "82.16"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.91%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.001"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.06%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.54"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.54%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.177"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.12%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.80"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.32%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "26.462.85"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +6.43%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.433"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.21%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.397"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.80%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.400"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +17.89%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "152.14"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.72%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.36"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.62%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.910.13"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.46%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "130.82"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.95%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.122"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.99%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.961"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.06%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08601"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.39%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.698"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.56%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "12.65"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.24%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.335"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.96%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02319"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.86%  "
$ws.Range("E39").Value = "  +3.47%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06185"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.42%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.432"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.71%  "
$ws.Range("E42").Value = "  -4.90%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.6159"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.52%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.000"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.02%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "14.08"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.04%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.832"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.20%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5944"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +5.12%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "128.22"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.44%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.032"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.95%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07160"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.43%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "76.85"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.39%  "
